# Generate Report for Handoff
# The source file "812a8b29-3a0a-4ecb-bf20-47442c426e62.md" has finished its
# translation round and is now ready to be handed off. Update the status /
# priority / timestamp columns on every sheet to reflect that.

$wb = $excel.ActiveWorkbook

$ws_overview = $wb.Worksheets.Item("Overview")
$ws_zhcn     = $wb.Worksheets.Item("zh-cn")
$ws_dede     = $wb.Worksheets.Item("de-de")

# --- Overview sheet: row 3 is the 812a8b29-...md file -------------------
$ws_overview.Range("E3").Value = "Ready for handoff"
$ws_overview.Range("F3").Value = "Ready for handoff"
$ws_overview.Range("G3").Value = "2016-09-01 10:16:15"

# --- zh-cn sheet: row 3 is the 812a8b29-...md file -----------------------
$ws_zhcn.Range("C3").Value = "Ready for handoff"
$ws_zhcn.Range("E3").Value = "mt"
$ws_zhcn.Range("H3").Value = "2016-09-01 10:16:04"

# --- de-de sheet: row 3 is the 812a8b29-...md file -----------------------
$ws_dede.Range("C3").Value = "Ready for handoff"
$ws_dede.Range("E3").Value = "mt"
$ws_dede.Range("H3").Value = "2016-09-01 10:16:15"

# --- Widen the columns that now hold the longer "Ready for handoff" /
#     timestamp text so nothing is truncated -----------------------------
$ws_overview.Range("E1:F1").ColumnWidth = 17.2159881591797
$ws_zhcn.Range("C1").ColumnWidth = 17.2159881591797
$ws_dede.Range("C1").ColumnWidth = 17.2159881591797
